$d = $word.ActiveDocument

# --- Change 1: "Objetivos" paragraph -----------------------------------
# Split the single run into two <w:t> runs joined by a manual line break
# right before the "2 Fornecer subsídios..." sentence.
$find1 = "1. Introduzir e discutir conceitos e técnicas estatísticas para controle e melhoria da qualidade de produtos fabricados e processos de fabricação;2 Fornecer subsídios para que o aluno tenha condições de utilizar essas técnicas e conceitos na sua vida profissional futura."
$repl1 = "1. Introduzir e discutir conceitos e técnicas estatísticas para controle e melhoria da qualidade de produtos fabricados e processos de fabricação;^l2 Fornecer subsídios para que o aluno tenha condições de utilizar essas técnicas e conceitos na sua vida profissional futura."

$d.Content.Find.Execute($find1, $false, $false, $false, $false, $false, $true, 1, $false, $repl1, 2)

# --- Change 2: "Programa" paragraph -------------------------------------
# Break the single run of text into one <w:t> per topic/sub-topic,
# separated by manual line breaks (<w:br/>), with an extra blank line
# (double break) between each top-level numbered section.
$find2 = "1. Sistemas de Medição.1.1. Planejamento do Sistema de Medição;1.2. Impacto da Variabilidade do Sistema de Medição no Produto;1.3. Sistemas de Medição por Atributos;1.4. Tendência e Linearidade;1.5. Análise de Repetitividade e Reprodutibilidade;2. Fundamentos do Controle Estatístico da Qualidade e do Processo.2.1. Importância do Controle Estatístico da Qualidade e do Processo;2.2. Naturezas das Variações;2.3. Causas Comuns e Causas Especiais de Variações;3. Gráficos de Controle por Variáveis3.1. Gráficos de Controle por Médias;3.2. Gráficos de Controle por Amplitude;3.3. Gráficos de Controle por Desvio Padrão;3.4. Análise de Desempenho dos Gráficos de Controle por Variáveis;4. Gráficos de Controle por Atributos4.1. Gráficos de Controle por Número de Não Conformidades;4.2. Gráficos de Controle por Fração Não Conforme;4.3. Gráficos de Controle por Número de Defeitos4.4. Gráficos de Controle por Não Conformidades por Amostra;5. Gráficos de Controle para Processos Auto-correlacionados5.1. Gráficos de Controle por Amplitude Móvel;5.2. Gráficos de Controle por Soma Acumulada (CUSUM).5.3. Gráficos de Controle por Média Móvel Ponderada Exponencialmente (EWMA)6. Analise de Capacidade do Processo6.1. Índices de Capacidade do Processo;6.2. Índices de Performance do Processo; 7. Inspeção da Qualidade7.1. Planos de Amostragem7.2. Inspeção para Aceitação;7.3. Inspeção Retificadora;8. Estudos de casos"

$repl2 = "1. Sistemas de Medição.^l" + `
  "1.1. Planejamento do Sistema de Medição;^l" + `
  "1.2. Impacto da Variabilidade do Sistema de Medição no Produto;^l" + `
  "1.3. Sistemas de Medição por Atributos;^l" + `
  "1.4. Tendência e Linearidade;^l" + `
  "1.5. Análise de Repetitividade e Reprodutibilidade;^l^l" + `
  "2. Fundamentos do Controle Estatístico da Qualidade e do Processo.^l" + `
  "2.1. Importância do Controle Estatístico da Qualidade e do Processo;^l" + `
  "2.2. Naturezas das Variações;^l" + `
  "2.3. Causas Comuns e Causas Especiais de Variações;^l^l" + `
  "3. Gráficos de Controle por Variáveis^l" + `
  "3.1. Gráficos de Controle por Médias;^l" + `
  "3.2. Gráficos de Controle por Amplitude;^l" + `
  "3.3. Gráficos de Controle por Desvio Padrão;^l" + `
  "3.4. Análise de Desempenho dos Gráficos de Controle por Variáveis;^l^l" + `
  "4. Gráficos de Controle por Atributos^l" + `
  "4.1. Gráficos de Controle por Número de Não Conformidades;^l" + `
  "4.2. Gráficos de Controle por Fração Não Conforme;^l" + `
  "4.3. Gráficos de Controle por Número de Defeitos^l" + `
  "4.4. Gráficos de Controle por Não Conformidades por Amostra;^l^l" + `
  "5. Gráficos de Controle para Processos Auto-correlacionados^l" + `
  "5.1. Gráficos de Controle por Amplitude Móvel;^l" + `
  "5.2. Gráficos de Controle por Soma Acumulada (CUSUM).^l" + `
  "5.3. Gráficos de Controle por Média Móvel Ponderada Exponencialmente (EWMA)^l^l" + `
  "6. Analise de Capacidade do Processo^l" + `
  "6.1. Índices de Capacidade do Processo;^l" + `
  "6.2. Índices de Performance do Processo; ^l^l" + `
  "7. Inspeção da Qualidade^l" + `
  "7.1. Planos de Amostragem^l" + `
  "7.2. Inspeção para Aceitação;^l" + `
  "7.3. Inspeção Retificadora;^l^l" + `
  "8. Estudos de casos"

$d.Content.Find.Execute($find2, $false, $false, $false, $false, $false, $true, 1, $false, $repl2, 2)
